$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.087760546121217
$ws.Range("C2").Value = 0.04825589860627133
$ws.Range("E2").Value = 0.04750515606672234
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002585252955586159
$ws.Range("I2").Value = 1.919155205181539
$ws.Range("K2").Value = 0.9120395085126347
$ws.Range("L2").Value = 0.2289908694571707
$ws.Range("M2").Value = 0.2669324270465125
$ws.Range("N2").Value = 3.368090532953836

$ws.Range("B3").Value = 1.055881317380909
$ws.Range("C3").Value = 0.04303695656498974
$ws.Range("E3").Value = 0.0475541415498828
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002588993524376928
$ws.Range("I3").Value = 1.923327444675095
$ws.Range("K3").Value = 0.8769177771089858
$ws.Range("L3").Value = 0.2264469966731184
$ws.Range("M3").Value = 0.2609855045765883
$ws.Range("N3").Value = 3.388522768829588

$ws.Range("B4").Value = 1.036902547631684
$ws.Range("C4").Value = 0.03985125068159334
$ws.Range("E4").Value = 0.04759677392551875
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.00259141213044032
$ws.Range("I4").Value = 1.926603453672207
$ws.Range("K4").Value = 0.855877347986933
$ws.Range("L4").Value = 0.224994802584682
$ws.Range("M4").Value = 0.2574755925714953
$ws.Range("N4").Value = 3.401907344600104

$ws.Range("B5").Value = 1.029318433428671
$ws.Range("C5").Value = 0.0385576726479826
$ws.Range("E5").Value = 0.04761731357747756
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002592428476196738
$ws.Range("I5").Value = 1.928118036153698
$ws.Range("K5").Value = 0.84743511696675
$ws.Range("L5").Value = 0.2244306541614236
$ws.Range("M5").Value = 0.2560809107140933
$ws.Range("N5").Value = 3.407572596812294

$ws.Range("B6").Value = 1.028068156521499
$ws.Range("C6").Value = 0.0383431515764272
$ws.Range("E6").Value = 0.04762091570789462
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002592599099436352
$ws.Range("I6").Value = 1.92838037727001
$ws.Range("K6").Value = 0.8460412594631634
$ws.Range("L6").Value = 0.224338647689116
$ws.Range("M6").Value = 0.2558514786291681
$ws.Range("N6").Value = 3.408526045238176

$ws.Range("B7").Value = 1.036799658429999
$ws.Range("C7").Value = 0.03983378643690116
$ws.Range("E7").Value = 0.04759703809597582
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002591425712743238
$ws.Range("I7").Value = 1.926623152769942
$ws.Range("K7").Value = 0.8557629589398914
$ws.Range("L7").Value = 0.2249870823528468
$ws.Range("M7").Value = 0.2574566390510995
$ws.Range("N7").Value = 3.401982894233797

$ws.Range("B8").Value = 1.076645170747952
$ws.Range("C8").Value = 0.04645246046665363
$ws.Range("E8").Value = 0.0475194461011279
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002586517465048233
$ws.Range("I8").Value = 1.920445539033842
$ws.Range("K8").Value = 0.8998206917717368
$ws.Range("L8").Value = 0.2280909815213192
$ws.Range("M8").Value = 0.2648525817058776
$ws.Range("N8").Value = 3.374961371293992

$ws.Range("B9").Value = 1.159500748984556
$ws.Range("C9").Value = 0.05958514405088522
$ws.Range("E9").Value = 0.04746651352291309
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002577855093571557
$ws.Range("I9").Value = 1.914000623620019
$ws.Range("K9").Value = 0.9903833977892589
$ws.Range("L9").Value = 0.2350476801351959
$ws.Range("M9").Value = 0.2804778945916766
$ws.Range("N9").Value = 3.328632966748316

$ws.Range("B10").Value = 1.223254214587485
$ws.Range("C10").Value = 0.06933480250528135
$ws.Range("E10").Value = 0.04748762136869544
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002572071575130529
$ws.Range("I10").Value = 1.91272681111581
$ws.Range("K10").Value = 1.059473802901806
$ws.Range("L10").Value = 0.2406888170787767
$ws.Range("M10").Value = 0.2926420275941553
$ws.Range("N10").Value = 3.298659533339759

$ws.Range("B11").Value = 1.252883768583189
$ws.Range("C11").Value = 0.07379379874066672
$ws.Range("E11").Value = 0.04751015061986408
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002569565291282631
$ws.Range("I11").Value = 1.912900185453871
$ws.Range("K11").Value = 1.091463094742124
$ws.Range("L11").Value = 0.2433702135859193
$ws.Range("M11").Value = 0.2983245569976418
$ws.Range("N11").Value = 3.285907143506009

$ws.Range("B12").Value = 1.264193911927066
$ws.Range("C12").Value = 0.0754858462797614
$ws.Range("E12").Value = 0.04752053189641003
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002568634054653
$ws.Range("I12").Value = 1.913074174925434
$ws.Range("K12").Value = 1.103657219062228
$ws.Range("L12").Value = 0.2444021389457163
$ws.Range("M12").Value = 0.3004977905412645
$ws.Range("N12").Value = 3.281205179267488

$ws.Range("B13").Value = 1.26175406725929
$ws.Range("C13").Value = 0.07512127537199831
$ws.Range("E13").Value = 0.04751821396857636
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002568833821146461
$ws.Range("I13").Value = 1.913031883595878
$ws.Range("K13").Value = 1.101027417287639
$ws.Range("L13").Value = 0.2441791602514627
$ws.Range("M13").Value = 0.300028795486206
$ws.Range("N13").Value = 3.282212176454109

$ws.Range("B14").Value = 1.253812456792588
$ws.Range("C14").Value = 0.07393293336319573
$ws.Range("E14").Value = 0.0475109676694121
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002569488320866881
$ws.Range("I14").Value = 1.912912327986334
$ws.Range("K14").Value = 1.092464699937722
$ws.Range("L14").Value = 0.2434547794112092
$ws.Range("M14").Value = 0.2985029219296109
$ws.Range("N14").Value = 3.285517762084865

$ws.Range("B15").Value = 1.24895971553056
$ws.Range("C15").Value = 0.07320550115547064
$ws.Range("E15").Value = 0.04750676975339019
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002569891541663822
$ws.Range("I15").Value = 1.91285320771695
$ws.Range("K15").Value = 1.087230268046056
$ws.Range("L15").Value = 0.2430132282014767
$ws.Range("M15").Value = 0.2975710637311977
$ws.Range("N15").Value = 3.287559084340202

$ws.Range("B16").Value = 1.221330468951805
$ws.Range("C16").Value = 0.06904388428601749
$ws.Range("E16").Value = 0.04748640838849205
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002572237867852488
$ws.Range("I16").Value = 1.912730636421074
$ws.Range("K16").Value = 1.057394489223356
$ws.Range("L16").Value = 0.2405158974357278
$ws.Range("M16").Value = 0.2922736555318721
$ws.Range("N16").Value = 3.299510707380534

$ws.Range("B17").Value = 1.204541467718144
$ws.Range("C17").Value = 0.06649703814743191
$ws.Range("E17").Value = 0.04747722190587567
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.0025737091301433
$ws.Range("I17").Value = 1.912848318189596
$ws.Range("K17").Value = 1.039234552375433
$ws.Range("L17").Value = 0.239013354397585
$ws.Range("M17").Value = 0.289061999455221
$ws.Range("N17").Value = 3.307068828726301

$ws.Range("B18").Value = 1.194943979858834
$ws.Range("C18").Value = 0.06503439661433674
$ws.Range("E18").Value = 0.04747315538503827
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002574567101065676
$ws.Range("I18").Value = 1.912986861363336
$ws.Range("K18").Value = 1.028842121068919
$ws.Range("L18").Value = 0.2381599776787056
$ws.Range("M18").Value = 0.2872287673749767
$ws.Range("N18").Value = 3.31149913181369

$ws.Range("B19").Value = 1.191704592296077
$ws.Range("C19").Value = 0.06453955261201827
$ws.Range("E19").Value = 0.04747198787808848
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002574859613993733
$ws.Range("I19").Value = 1.913045936887251
$ws.Range("K19").Value = 1.025332471400958
$ws.Range("L19").Value = 0.2378729027396673
$ws.Range("M19").Value = 0.2866104770799964
$ws.Range("N19").Value = 3.313013419117027

$ws.Range("B20").Value = 1.206322570227826
$ws.Range("C20").Value = 0.06676792205166748
$ws.Range("E20").Value = 0.04747807390191028
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002573551297447735
$ws.Range("I20").Value = 1.912828456890416
$ws.Range("K20").Value = 1.041162256160135
$ws.Range("L20").Value = 0.2391721803306268
$ws.Range("M20").Value = 0.2894024342472434
$ws.Range("N20").Value = 3.306255655039251

$ws.Range("B21").Value = 1.256142657279213
$ws.Range("C21").Value = 0.07428188185028262
$ws.Range("E21").Value = 0.04751304594756611
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002569295595062189
$ws.Range("I21").Value = 1.912944503502686
$ws.Range("K21").Value = 1.094977593639527
$ws.Range("L21").Value = 0.2436670989974346
$ws.Range("M21").Value = 0.2989505279968157
$ws.Range("N21").Value = 3.284543380986833

$ws.Range("B22").Value = 1.289227895961631
$ws.Range("C22").Value = 0.07921326536960294
$ws.Range("E22").Value = 0.04754668093620218
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002566618183555793
$ws.Range("I22").Value = 1.913651852840005
$ws.Range("K22").Value = 1.130618086817179
$ws.Range("L22").Value = 0.2467011751859616
$ws.Range("M22").Value = 0.3053153759216301
$ws.Range("N22").Value = 3.271094005049719

$ws.Range("B23").Value = 1.271521715232211
$ws.Range("C23").Value = 0.07657937855771024
$ws.Range("E23").Value = 0.04752774600327569
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002568037687380108
$ws.Range("I23").Value = 1.913216517401708
$ws.Range("K23").Value = 1.111553169318483
$ws.Range("L23").Value = 0.2450730217001365
$ws.Range("M23").Value = 0.3019069485807861
$ws.Range("N23").Value = 3.278204347751924

$ws.Range("B24").Value = 1.205517163279382
$ws.Range("C24").Value = 0.06664545054401572
$ws.Range("E24").Value = 0.04747768492981486
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002573622615903011
$ws.Range("I24").Value = 1.912837215370622
$ws.Range("K24").Value = 1.040290591892301
$ws.Range("L24").Value = 0.2391003425581317
$ws.Range("M24").Value = 0.2892484825834742
$ws.Range("N24").Value = 3.306623025960121

$ws.Range("B25").Value = 1.136580675284364
$ws.Range("C25").Value = 0.05601519952408296
$ws.Range("E25").Value = 0.04747026165365575
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002580096070573106
$ws.Range("I25").Value = 1.915136686710014
$ws.Range("K25").Value = 0.9654363420362415
$ws.Range("L25").Value = 0.2330726063093422
$ws.Range("M25").Value = 0.276130701227487
$ws.Range("N25").Value = 3.340452608446697
